$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Aula" (room) numbers for the two course rows. ---
# The original cells store these as text (shared strings), not numbers,
# so force a Text number format before writing the new value; otherwise
# the numeric-looking string would be auto-coerced into a real number.
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "171"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "350"

# --- Recolor the header rows: light green (CDFFC4) -> light yellow (FFF8C4) ---
# Excel COM's Color is a BGR-ordered long (0x00BBGGRR), so build it from
# the target RGB hex with the byte order reversed.
# (Target the title cell and the column-header row separately rather than
# the whole merged A1:E1 range, so we don't spuriously materialize empty
# cells B1:E1 that the merge had left out of the sheet data.)
$newColor = [System.Convert]::ToInt32("C4F8FF", 16)
$ws.Range("A1").Interior.Color = $newColor
$ws.Range("A2:E2").Interior.Color = $newColor
